$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell, far away from the used range, used to mint a clean shared
# string (no quotePrefix/number-format style) holding the new numeric-looking
# text value. A direct $ws.Range("A2").Value = "1410001" assignment would be
# parsed as a number (losing the shared-string text type), and prefixing with
# an apostrophe stamps a quotePrefix style onto the cell that the target
# workbook doesn't have - so we go through a text formula + paste-values
# instead, which preserves plain text typing with no style changes.
$scratch = $ws.Range("Z100")
$scratch.Formula = '="1410001"'
$scratch.Copy()
$ws.Range("A2").PasteSpecial(-4163)
$ws.Range("A3").PasteSpecial(-4163)
$scratch.ClearContents()
